$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new blank column before column A. This shifts the existing
#    columns A:K to B:L, carrying their values AND styles with them
#    (old A2:A20 "segment name" cells -> B2:B20, old B1:K1 headers -> C1:L1,
#    old B2:K20 data values -> C2:L20).
$ws.Columns.Item(1).Insert()

# 2) New column B1 needs to become the "segments" header (bold / bordered /
#    centered, same style as the other header cells). Copy formatting from
#    the now-shifted header cell C1 (style s=1) into B1, then set its text.
$ws.Range("C1").Copy($ws.Range("B1"))
$ws.Range("B1").Value = "segments"

# 3) New column A2:A20 needs to hold the numeric segment index (0-18) using
#    the same styled format (s=1) that the segment-name cells currently
#    sitting in B2:B20 have (inherited from the pre-insert column A). Copy
#    that formatting into A2:A20 before overwriting B with plain values.
$xlDown = -4121
$lastRow = $ws.Cells(1, 2).End($xlDown).Row
For ($r = 2; $r -le $lastRow; $r++) {
    $ws.Range("B" + $r).Copy($ws.Range("A" + $r))
    $ws.Range("A" + $r).Value = $r - 2
}

# 4) The segment-name text cells now living in column B should end up with
#    plain (unstyled) formatting, matching the target layout.
$ws.Range("B2:B" + $lastRow).Style = "Normal"
